$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Round the values in row 5 to 2 decimal places (custom accuracy)
$row5 = @{
    "B5" = 20.03
    "C5" = 15.05
    "D5" = 0.85
    "E5" = 43.78
    "F5" = 35.84
    "G5" = 15.45
    "H5" = 59.22
    "I5" = 24.3
    "J5" = 11.21
    "K5" = 16.07
    "L5" = 17.53
    "M5" = 18.75
    "N5" = 5.36
    "O5" = 15.79
    "P5" = 22.31
    "Q5" = 13.23
    "R5" = 0.22
    "S5" = 0.77
    "T5" = 233.02
    "U5" = 44.1
    "V5" = 14.58
    "W5" = 29.59
    "X5" = 15.83
    "Y5" = 2.08
    "Z5" = 29.74
    "AA5" = 12.72
    "AB5" = 11.54
    "AC5" = 13.59
    "AD5" = 18.73
    "AE5" = 0.34
    "AF5" = 53.87
    "AG5" = 8.460000000000001
    "AH5" = 18.06
}

foreach ($addr in $row5.Keys) {
    $ws.Range($addr).Value = $row5[$addr]
}

# 2) Remove the last data row (row 6) - "데이터 1000개" trimming
$ws.Rows.Item(6).Delete()

# 3) Shrink columns C and E to width 7 (custom accuracy column formatting)
#    Excel's ColumnWidth property (character units) gets converted on save to the
#    OOXML stored "width" attribute with an offset of 5/6; using 6.166666666666667
#    here round-trips to a stored width of exactly 7.
$ws.Columns.Item(3).ColumnWidth = 6.166666666666667
$ws.Columns.Item(5).ColumnWidth = 6.166666666666667
